$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 282, shifting existing rows 282-309 down to 283-310
$ws.Rows.Item(282).Insert()

# Populate the newly inserted row 282 with the new weekly record
$ws.Cells.Item(282, 1).Value = 9
$ws.Cells.Item(282, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(282, 3).Value = "Metropolitana"
$ws.Cells.Item(282, 4).Value = 44748
$ws.Cells.Item(282, 5).Value = 13
$ws.Cells.Item(282, 6).Value = 100112021
$ws.Cells.Item(282, 7).Value = "Ají"
$ws.Cells.Item(282, 8).Value = "Americana (o)"
$ws.Cells.Item(282, 9).Value = "Primera"
$ws.Cells.Item(282, 10).Value = 25
$ws.Cells.Item(282, 11).Value = 35000
$ws.Cells.Item(282, 12).Value = 36000
$ws.Cells.Item(282, 13).Value = 35520
$ws.Cells.Item(282, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(282, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(282, 16).Value = 1421
$ws.Cells.Item(282, 17).Value = 25
$ws.Cells.Item(282, 18).Value = "Hortaliza"
